# Resize the 4 columns of the "User Stories" table.
# Word's Column.Width / Cell.Width are expressed in points, while the
# underlying OOXML <w:tcW>/<w:gridCol> values are in twentieths of a point
# (twips), so divide the target twip values by 20 before assigning.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$table.Columns.Item(1).Width = 1129 / 20
$table.Columns.Item(2).Width = 2410 / 20
$table.Columns.Item(3).Width = 2693 / 20
$table.Columns.Item(4).Width = 2784 / 20
